# ADDITION graph for 'Góc che khuất' tab pane.
#
# Adds a new worksheet "Sheet2" right after "Sheet1", fills it with the
# B1:C11 data table used to drive the "Góc che khuất" chart (B: 0..5 step
# 0.5, C: 1..11), makes it the active/selected tab, and leaves the
# selection parked on G9 (mirrors the author's on-screen selection when
# the sheet was added). Sheet1's own selection (D9) is left untouched.

$wb = $excel.ActiveWorkbook

# Add the new sheet immediately after the existing "Sheet1" so it becomes
# tab index 2 / activeTab=1, matching the workbook-level bookViews change.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# B1:C11 data table: column B is 0, 0.5, 1, ... 5 (step 0.5); column C is
# the row index 1..11.
$bValues = @(0, 0.5, 1, 1.5, 2, 2.5, 3, 3.5, 4, 4.5, 5)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 2).Value = $bValues[$i]
    $ws2.Cells.Item($row, 3).Value = $row
}

# New sheet becomes the active tab, selection sits at G9.
[void]$ws2.Range("G9").Select()
